# The document currently ends with:
#   ... <drawing UC11.png> <bookmarkStart _GoBack/><bookmarkEnd _GoBack/> </w:p>
#   <w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>   (x4, empty trailing paragraphs)
#
# Target layout:
#   ... <drawing UC11.png> </w:p>
#   <w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="24"/>
#                       <w:lang w:val="en-US"/></w:rPr></w:pPr>
#        <bookmarkStart _GoBack/><bookmarkEnd _GoBack/></w:p>
#   (the four empty trailing paragraphs are removed)

$d = $word.ActiveDocument

# 1) Locate the (hidden) _GoBack bookmark that currently sits at the end of the
#    paragraph holding the UC11 picture, remember its position, then remove it.
$bookmark = $d.Bookmarks.Item("_GoBack")
$splitPos = $bookmark.Start
$bookmark.Delete()

# 2) Split that paragraph right where the bookmark used to be: inserting a
#    paragraph mark there closes the picture's paragraph and opens a new one
#    that inherits the same run formatting (bold, size 28/24, en-US) from the
#    picture paragraph's mark - matching the new <w:pPr><w:rPr> in the diff.
$d.Range($splitPos, $splitPos).InsertAfter("`r")
$newParaPos = $splitPos + 1

# 3) Re-create the _GoBack bookmark inside that brand-new paragraph. A
#    collapsed (zero-length) range confuses Bookmarks.Add, so bracket the
#    insertion point with a throw-away character, bookmark the 1-char range,
#    then delete the character again, leaving a correctly placed bookmark.
$d.Range($newParaPos, $newParaPos).InsertAfter("X")
$markerRange = $d.Range($newParaPos, $newParaPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null
$d.Range($newParaPos, $newParaPos + 1).Delete()

# 4) Remove the four now-superfluous empty trailing paragraphs (the ones whose
#    paragraph mark only carries <w:lang w:val="en-US"/>) that used to follow
#    the picture paragraph - the new bookmark paragraph is the last one now.
#    Resolve the 1-based Paragraphs index of the freshly-created bookmark
#    paragraph by matching its start offset (Paragraphs.Item only accepts an
#    index, not a character position).
$newBookmarkIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $newParaPos) {
        $newBookmarkIndex = $i
        break
    }
}

$firstEmptyPara = $d.Paragraphs.Item($newBookmarkIndex + 1)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Range($firstEmptyPara.Range.Start, $lastPara.Range.End).Delete()
